# "Creacion de seccion e-commerce"
# The header row is being reworked for the new e-commerce section:
#   - Multimedia_1 (column E) is replaced by a new "Imagen_principal" header.
#   - Multimedia_2/3/4 keep their text and simply shift to fill column E's
#     old neighbours (handled automatically by the shared-string table once
#     the E1 text changes).
#   - The active selection on the sheet moves from F8 to G4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Imagen_principal"

$ws.Range("G4").Select()
